# Updated cryptos list values to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''69.944.88'
$ws.Range('E2').Value = '''  -0.32%  '
$ws.Range('D3').Value = '''3.539.32'
$ws.Range('E3').Value = '''  -0.47%  '
$ws.Range('E4').Value = '''  -0.14%  '
$ws.Range('D5').Value = '''603.30'
$ws.Range('E5').Value = '''  -2.29%  '
$ws.Range('D6').Value = '''196.63'
$ws.Range('E6').Value = '''  +5.77%  '
$ws.Range('D7').Value = '''0.626'
$ws.Range('E7').Value = '''  -0.25%  '
$ws.Range('E8').Value = '''  -0.09%  '
$ws.Range('E9').Value = '''  -2.73%  '
$ws.Range('E10').Value = '''  +0.22%  '
$ws.Range('D11').Value = '''54.08'
$ws.Range('E11').Value = '''  +0.32%  '
$ws.Range('D12').Value = '''0.0000303'
$ws.Range('E12').Value = '''  -2.42%  '
$ws.Range('D13').Value = '''9.56'
$ws.Range('E13').Value = '''  +0.21%  '
$ws.Range('D14').Value = '''4.089.87'
$ws.Range('E14').Value = '''  -0.71%  '
$ws.Range('D15').Value = '''600.07'
$ws.Range('E15').Value = '''  -4.58%  '
$ws.Range('D16').Value = '''70.107.73'
$ws.Range('E16').Value = '''  -0.17%  '
$ws.Range('D17').Value = '''19.16'
$ws.Range('E17').Value = '''  +1.11%  '
$ws.Range('D18').Value = '''12.68'
$ws.Range('E18').Value = '''  -2.92%  '
$ws.Range('D19').Value = '''3.524.02'
$ws.Range('E19').Value = '''  -0.88%  '
$ws.Range('E20').Value = '''  +0.54%  '
$ws.Range('D21').Value = '''0.998'
$ws.Range('E21').Value = '''  +0.07%  '
$ws.Range('D22').Value = '''18.12'
$ws.Range('E22').Value = '''  +2.86%  '
$ws.Range('D23').Value = '''5.31'
$ws.Range('E23').Value = '''  +7.18%  '
$ws.Range('D24').Value = '''103.41'
$ws.Range('E24').Value = '''  +0.28%  '
$ws.Range('D25').Value = '''4.62'
$ws.Range('E25').Value = '''  -2.45%  '
$ws.Range('D26').Value = '''3.11'
$ws.Range('E26').Value = '''  +2.82%  '
$ws.Range('D27').Value = '''10.95'
$ws.Range('E27').Value = '''  -1.19%  '
$ws.Range('E28').Value = '''  +1.27%  '
$ws.Range('D29').Value = '''33.58'
$ws.Range('E29').Value = '''  -2.58%  '
$ws.Range('D30').Value = '''4.48'
$ws.Range('E30').Value = '''  +19.59%  '
$ws.Range('D31').Value = '''7.12'
$ws.Range('E31').Value = '''  +0.50%  '
$ws.Range('D32').Value = '''12.76'
$ws.Range('E32').Value = '''  +3.71%  '
$ws.Range('E33').Value = '''  +0.98%  '
$ws.Range('D34').Value = '''63.52'
$ws.Range('E34').Value = '''  -0.50%  '
$ws.Range('D35').Value = '''0.0₃0829'
$ws.Range('E35').Value = '''  +5.77%  '
$ws.Range('D36').Value = '''3.744.43'
$ws.Range('E36').Value = '''  +4.59%  '
$ws.Range('D37').Value = '''3.10'
$ws.Range('E37').Value = '''  -5.31%  '
$ws.Range('E38').Value = '''  +0.01%  '
$ws.Range('D39').Value = '''0.395'
$ws.Range('E39').Value = '''  -2.04%  '
$ws.Range('D40').Value = '''3.61'
$ws.Range('E40').Value = '''  +1.92%  '
$ws.Range('D41').Value = '''36.88'
$ws.Range('E41').Value = '''  -1.12%  '
$ws.Range('D42').Value = '''498.29'
$ws.Range('E42').Value = '''  -6.29%  '
$ws.Range('E43').Value = '''  -1.08%  '
$ws.Range('E44').Value = '''  -1.04%  '
$ws.Range('E45').Value = '''  -1.94%  '
$ws.Range('B46').Value = '''ThetaToken'
$ws.Range('C46').Value = '''https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = '''2.83'
$ws.Range('E46').Value = '''  -4.16%  '
$ws.Range('B47').Value = '''ApeXProtocol'
$ws.Range('C47').Value = '''https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''3.32'
$ws.Range('E47').Value = '''  -1.83%  '
$ws.Range('E48').Value = '''  +0.29%  '
$ws.Range('D49').Value = '''8.68'
$ws.Range('E49').Value = '''  -5.58%  '
$ws.Range('E50').Value = '''  +1.10%  '
$ws.Range('E51').Value = '''  +12.10%  '
